$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the old "2 pin female header" row entirely (old row 7).
#    This shifts the old row 8 ("30k surface mount resistors 0805")
#    up to become the new row 7, and drops the orphaned shared
#    strings / the pcb-sockets hyperlink relationship along with it.
# ------------------------------------------------------------------
$ws.Rows.Item(7).Delete()

# ------------------------------------------------------------------
# 2. Insert a new "Place Labels" column before the Comments column,
#    which is now column D - the old D (Comments) becomes E.
# ------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# ------------------------------------------------------------------
# 3. Populate the new "Place Labels" column, in row order, and fix
#    up the ADC part description. Order matches the commit so the
#    shared-string table is appended in the same sequence.
# ------------------------------------------------------------------
$ws.Range("D1").Value = "Place Labels"
$ws.Range("D2").Value = "C1, C2, ADC1Cap, ADC2Cap, Cref, Csr, C3, C4, C5, C6, C7, C8, ADCap1, ADCap2, OpCap"
$ws.Range("D3").Value = "Rr1, R2, Rr2"
$ws.Range("D4").Value = "R1,R3, R4, R5, R7, R6, R8, R9,R10,R11, R12, R13"
$ws.Range("A5").Value = "ADS1015IDGSR"
$ws.Range("D5").Value = "ADS1015-1, ADS1015-2"
$ws.Range("D6").Value = "AD623-1, AD623-2"
$ws.Range("D7").Value = "Rsr1, Rsr2"

# ------------------------------------------------------------------
# 4. Column widths: D (Place Labels) = 18, E (Comments) = 14.
#    ColumnWidth is specified in characters and Excel adds a fixed
#    padding offset (~0.8333 for this workbook's font) before it is
#    written out as the raw "width" attribute, so compensate for it.
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666

# ------------------------------------------------------------------
# 5. Hyperlinks: the engine's Hyperlinks collection does not track
#    row/column shifts, so its entries are now stale. Rebuild it
#    from scratch. Cell formatting (wrap/no-wrap "Hyperlink" style)
#    gets reset to a fresh style by Hyperlinks.Add, so stash/restore
#    the original formatting for every affected cell via PasteSpecial.
# ------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("Z1"))
$ws.Range("B3").Copy($ws.Range("Z2"))
$ws.Range("B4").Copy($ws.Range("Z3"))
$ws.Range("B7").Copy($ws.Range("Z4"))
$ws.Range("B5").Copy($ws.Range("Z5"))

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "http://china.rs-online.com/web/p/ceramic-multilayer-capacitors/6911161/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6789667/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6791569/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6792039/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://china.rs-online.com/web/p/general-purpose-adcs/7094550/") | Out-Null

$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").Copy()
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy()
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy()
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("Z5").Copy()
$ws.Range("B5").PasteSpecial(-4122) | Out-Null

$ws.Range("Z1:Z5").Clear()
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 6. Selection, matching the target workbook view.
# ------------------------------------------------------------------
$ws.Range("G11").Select()

$wb.Save()
